# Apply the "Add files via upload" edit to cryptoPaxgoldPolygon.xlsx
#   - D2: 10 -> 20              (combien de dollars j'ajoute a chaque achat)
#   - A4: 4000 -> 4050          (new latest price sample)
#   - B4: 2E-11 -> 4            (new latest quantity sample)
#   - Rows 5-8 (the old price-history rows) are wiped out entirely so the
#     summary formulas in I2:M2 recompute purely from row 4
#   - Selection moves to C12, matching the saved view state

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "combien de dollars j'ajoute a chaque achat" input
$ws.Range("D2").Value = 20

# Refresh the latest price/quantity sample on row 4
$ws.Range("A4").Value = 4050
$ws.Range("B4").Value = 4

# Wipe out the old historical rows (5-8) completely - values, text and
# formatting - so the sheet ends up with only the single data row (4)
$ws.Range("A5:D8").Clear()

# Restore the saved selection
$ws.Range("C12").Select()
